# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (AD1:AF1) ---------------------------------------------------
# Copy the formatting of the existing header cell (bold, bordered, centered
# style index 1) onto the new header cells before setting their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (AD2:AF54) -------------------------------------------------
# Every player row gets the same season record: 57 wins, 105 losses, 0 ties.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 57
    $ws.Cells.Item($r, 31).Value = 105
    $ws.Cells.Item($r, 32).Value = 0
}
